$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update B26: calibration materials reference bumped from 37 -> 42, with date note added
$ws.Range("B26").Value = "using materials 42 (23 August 2020 calibration) material right now"

# New note added in B27, explaining a redo that couldn't be completed
$ws.Range("B27").Value = "couldn't redo this, need to fix! So July 30 figures there"

# Update the active selection to A12 (cursor moved there after the edits)
$ws.Range("A12").Select()
